# "Archive / small fixed"
# Appends one new archive-log row (row 46) to the tracker sheet:
#   A46 = 2019-03-04 (serial 43528), formatted as a short date
#   B46 = "记录存档逻辑"  ("Archive logic record")
# and moves the active selection to B47 (just past the new row),
# matching the author's next-entry cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date cell - same short-date number format (built-in id 14) used
# throughout column A.
$ws.Range("A46").Value = 43528
$ws.Range("A46").NumberFormat = "mm-dd-yy"

# New note cell - becomes a new shared-string entry.
$ws.Range("B46").Value = "记录存档逻辑"

# Leave the selection where the author would continue typing next.
$ws.Range("B47").Select()
